$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (I1, J1) ---
# Copy formatting from the existing header cell H1 (bold, centered, bordered)
# so the new header cells I1/J1 share the same style index as the rest
# of row 1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (I2:J69) ---
$data = @(
    @(2,10,10),
    @(3,7,7),
    @(4,7,7),
    @(5,7,8),
    @(6,8,8),
    @(7,8,8),
    @(8,7,7),
    @(9,8,8),
    @(10,7,8),
    @(11,9,9),
    @(12,9,9),
    @(13,8,9),
    @(14,9,9),
    @(15,10,10),
    @(16,10,10),
    @(17,7,7),
    @(18,6,6),
    @(19,6,6),
    @(20,7,7),
    @(21,8,8),
    @(22,6,6),
    @(23,9,9),
    @(24,7,7),
    @(25,6,6),
    @(26,9,9),
    @(27,7,8),
    @(28,9,9),
    @(29,7,7),
    @(30,7,7),
    @(31,8,8),
    @(32,7,8),
    @(33,6,6),
    @(34,6,6),
    @(35,8,8),
    @(36,5,7),
    @(37,7,7),
    @(38,5,6),
    @(39,9,9),
    @(40,9,9),
    @(41,6,7),
    @(42,4,5),
    @(43,8,8),
    @(44,6,7),
    @(45,8,8),
    @(46,6,6),
    @(47,7,8),
    @(48,7,8),
    @(49,7,7),
    @(50,6,6),
    @(51,8,8),
    @(52,7,7),
    @(53,8,8),
    @(54,7,8),
    @(55,8,8),
    @(56,6,6),
    @(57,6,7),
    @(58,8,9),
    @(59,5,6),
    @(60,7,8),
    @(61,8,8),
    @(62,7,8),
    @(63,7,7),
    @(64,5,5),
    @(65,5,6),
    @(66,6,7),
    @(67,6,7),
    @(68,4,5),
    @(69,6,6)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "I0/IF columns added"
